$wb = $excel.ActiveWorkbook

# Sheet2 = "Manila Philippines": row 4 values (Commit/Forecast, Professional Voluntary Turnover)
$ws2 = $wb.Worksheets.Item("Manila Philippines")
$ws2.Range("M4").Value = 0
$ws2.Range("N4").Value = 0
$ws2.Range("O4").Value = 0.0595166666666667
$ws2.Range("P4").Value = 0.0595166666666667
$ws2.Range("Q4").Value = 0.0595166666666667
$ws2.Range("R4").Value = 0.17855
$ws2.Range("S4").Value = 0.0595166666666667
$ws2.Range("T4").Value = 0.0595166666666667
$ws2.Range("U4").Value = 0.0595166666666667
$ws2.Range("V4").Value = 0.17855
$ws2.Range("W4").Value = 0.7142

# Sheet3 = "Milwaukee Pmc Hq Wisconsin": row 7, clear M7:N7
$ws3 = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws3.Range("M7:N7").ClearContents()

# Sheet4 = "Milwaukee Wisconsin": row 5, clear M5:N5
$ws4 = $wb.Worksheets.Item("Milwaukee Wisconsin")
$ws4.Range("M5:N5").ClearContents()

# Sheet5 = "South Beloit Gardner St Illino": row 7, clear M7:N7
$ws5 = $wb.Worksheets.Item("South Beloit Gardner St Illino")
$ws5.Range("M7:N7").ClearContents()

# Sheet7 = "Rock Road Radford Virginia": row2 clear M2, row3 clear M3:N3
$ws7 = $wb.Worksheets.Item("Rock Road Radford Virginia")
$ws7.Range("M2").ClearContents()
$ws7.Range("M3:N3").ClearContents()
